$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New batch status rows continuing the existing pattern (rows 6-10)
$data = @(
    @(200, 250, "processed", 0.98),
    @(250, 300, "processed", 1),
    @(300, 350, "processed", 0.96),
    @(350, 400, "processed", 0.9399999999999999),
    @(400, 450, "processed", 0.88)
)

$startRow = 6
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $values = $data[$i]
    $ws.Cells.Item($row, 1).Value = $values[0]
    $ws.Cells.Item($row, 2).Value = $values[1]
    $ws.Cells.Item($row, 3).Value = $values[2]
    $ws.Cells.Item($row, 4).Value = $values[3]
}
